$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 : "1." (was "3.")
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "1."

# ---------------------------------------------------------------------------
# Row 3 : "1.1." child of "1."
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "1."

# ---------------------------------------------------------------------------
# Row 5 : description text updated
# ---------------------------------------------------------------------------
$ws.Range("F5").Value = "Propuesta de mejoras resultado de reuniones con gerencia y directores técnicos"

# ---------------------------------------------------------------------------
# Row 6 : "1.1.2." (was "1.1.2.1" style numbering shift)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "1.1.2."

# ---------------------------------------------------------------------------
# Row 7 : child of row 6, plus renamed name/description
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "1.1.2.1"
$ws.Range("B7").Value = "1.1.2."
$ws.Range("D7").Value = "Presentación implementación métodologia Last Planner System"
$ws.Range("F7").Value = "Presentación implementación Last Planner System"

# ---------------------------------------------------------------------------
# Row 9 : child of row 8
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "1.1.3.1."
$ws.Range("D9").Value = "Presentación implementación funcionalidad grilla Last Planner System"
$ws.Range("F9").Value = "Socialización e implementación de avances por grilla"

# ---------------------------------------------------------------------------
# Row 10 : "1.2." child of "1."
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "1.2."
$ws.Range("B10").Value = "1."

# ---------------------------------------------------------------------------
# Row 11 : "1.3." child of "1."
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "1.3."
$ws.Range("B11").Value = "1."

# ---------------------------------------------------------------------------
# Row 12 : "2."
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "2."

# ---------------------------------------------------------------------------
# Row 13 : now "2.1." Planeación Lean (used to be "3." Desarrollo Software)
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "2.1."
$ws.Range("B13").Value = "2."
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "Planeación Lean"

# ---------------------------------------------------------------------------
# Row 14 : new content "2.2." Auditoria (previously a fully blank row)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "2.2."
$ws.Range("B14").Value = "2."
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = "Auditoria"
$ws.Range("E14").Value = "Carpeta"

# ---------------------------------------------------------------------------
# Row 15 : new content "2.2.1." Excel de informe de Auditoria (blank before)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "2.2.1."
$ws.Range("B15").Value = "2.2."
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = "Excel de informe de Auditoria"
$ws.Range("E15").Value = "Xlsx"
$ws.Range("F15").Value = "Archivo de registro de auditorias de seguimiento"

# Hyperlink first (so the cell gets the plain hyperlink style), then copy the
# exact visual format from an existing hyperlink cell (G7) on top of it so
# the resulting style matches the rest of the table, and finally restore the
# correct display text/value.
$ws.Hyperlinks.Add($ws.Range("G15"), "https://docs.google.com/spreadsheets/d/1PvlOcqy-B7uOcPeKvaGO18cssIEnb6UIXeNBRuVQpiE/edit?usp=sharing") | Out-Null
$ws.Range("G7").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = "https://docs.google.com/spreadsheets/d/1PvlOcqy-B7uOcPeKvaGO18cssIEnb6UIXeNBRuVQpiE/edit?usp=sharing"

# ---------------------------------------------------------------------------
# Row 16 : "3." Desarrollo Software (used to be fully blank)
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "3."
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = "Desarrollo Software"
$ws.Range("E16").Value = "Carpeta"

# ---------------------------------------------------------------------------
# Extend formatting down to the three new trailing blank rows (21-23) so
# they carry the same shaded style as the rest of the table body.
# ---------------------------------------------------------------------------
$ws.Range("A20:H20").Copy()
$ws.Range("A21:H23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Resize the worksheet table / autofilter to cover the new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H23"))

# ---------------------------------------------------------------------------
# Restore the selection left behind by the edit.
# ---------------------------------------------------------------------------
$ws.Range("F12").Select()
